# "added this last report 09-01-25"
# Route Cost RSO.xlsx - update route costs for the latest report and refresh the report date.

$wb = $excel.ActiveWorkbook
$wsRoute = $wb.Worksheets.Item("Route")
$wsMobil = $wb.Worksheets.Item("Mobil")

# --- Route sheet: report date (L3) becomes a real date value instead of text ---
$wsRoute.Range("L3").Value = 45901

# --- Route sheet: updated daily route-cost figures for RSO 02 / RSO 03 / RSO 04 ---
$wsRoute.Range("D7").Value = 150
$wsRoute.Range("D9").Value = 150
$wsRoute.Range("D10").Value = 200

# --- Mobil sheet: bump the report date shown at the top ---
$wsMobil.Range("G1").Value = 45658

# --- Mobil sheet: move the on-sheet selection, without disturbing the active tab ---
$wsMobil.Range("G2:G3").Select()
$wsRoute.Activate()
